# Fruta / hortaliza, semanal
#
# Two new weekly price rows are inserted into the daily-logic subset sheet
# for "Vega Modelo de Temuco - Plátano", pushing the existing data for
# rows 995..1061 down to rows 997..1063 (dimension grows from T1061 to
# T1063). The two newly inserted rows (995 and 996) carry a new report
# date (45021) for "Barraganete/Verde" and "Sin especificar/Pintón"
# respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 995-1061 down by two rows, leaving two blank rows
# at 995:996 for the new data.
$ws.Rows("995:996").Insert()

# --- New row 995: Barraganete / Verde ---------------------------------
$ws.Cells.Item(995, 1).Value = 10
$ws.Cells.Item(995, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(995, 3).Value = "La Araucanía"
$ws.Cells.Item(995, 4).Value = 45021
$ws.Cells.Item(995, 5).Value = 9
$ws.Cells.Item(995, 6).Value = "Fruta"
$ws.Cells.Item(995, 7).Value = 100108
$ws.Cells.Item(995, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(995, 9).Value = 100108006
$ws.Cells.Item(995, 10).Value = "Plátano"
$ws.Cells.Item(995, 11).Value = "Barraganete"
$ws.Cells.Item(995, 12).Value = "Verde"
$ws.Cells.Item(995, 13).Value = 80
$ws.Cells.Item(995, 14).Value = 36000
$ws.Cells.Item(995, 15).Value = 36000
$ws.Cells.Item(995, 16).Value = 36000
$ws.Cells.Item(995, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(995, 18).Value = "Ecuador"
$ws.Cells.Item(995, 19).Value = 1800
$ws.Cells.Item(995, 20).Value = 20

# --- New row 996: Sin especificar / Pintón -----------------------------
$ws.Cells.Item(996, 1).Value = 10
$ws.Cells.Item(996, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(996, 3).Value = "La Araucanía"
$ws.Cells.Item(996, 4).Value = 45021
$ws.Cells.Item(996, 5).Value = 9
$ws.Cells.Item(996, 6).Value = "Fruta"
$ws.Cells.Item(996, 7).Value = 100108
$ws.Cells.Item(996, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(996, 9).Value = 100108006
$ws.Cells.Item(996, 10).Value = "Plátano"
$ws.Cells.Item(996, 11).Value = "Sin especificar"
$ws.Cells.Item(996, 12).Value = "Pintón"
$ws.Cells.Item(996, 13).Value = 800
$ws.Cells.Item(996, 14).Value = 25000
$ws.Cells.Item(996, 15).Value = 25000
$ws.Cells.Item(996, 16).Value = 25000
$ws.Cells.Item(996, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(996, 18).Value = "Ecuador"
$ws.Cells.Item(996, 19).Value = 1250
$ws.Cells.Item(996, 20).Value = 20
